$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A397").Value = "test"
